$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.28"
$ws.Range("E2").Value = "'1.05%"
$ws.Range("D3").Value = "'5.507"
$ws.Range("E3").Value = "'0.36%"
$ws.Range("D4").Value = "'0.08013"
$ws.Range("E4").Value = "'-0.32%"
$ws.Range("D5").Value = "'1.977"
$ws.Range("E5").Value = "'4.47%"
$ws.Range("D6").Value = "'4.335"
$ws.Range("E6").Value = "'0.87%"
$ws.Range("E7").Value = "'-2.67%"
$ws.Range("D8").Value = "'0.9487"
$ws.Range("E8").Value = "'0.92%"
$ws.Range("D9").Value = "'0.1119"
$ws.Range("E9").Value = "'-4.73%"
$ws.Range("D10").Value = "'0.1859"
$ws.Range("E10").Value = "'-0.59%"
$ws.Range("D11").Value = "'10.62"
$ws.Range("E11").Value = "'24.15%"
$ws.Range("D12").Value = "'0.09852"
$ws.Range("E12").Value = "'-0.78%"
$ws.Range("D13").Value = "'0.04676"
$ws.Range("E13").Value = "'11.21%"
$ws.Range("D14").Value = "'0.1066"
$ws.Range("E14").Value = "'-0.02%"
$ws.Range("D15").Value = "'0.001265"
$ws.Range("E15").Value = "'-0.54%"
$ws.Range("D16").Value = "'0.04066"
$ws.Range("E16").Value = "'-4.35%"
$ws.Range("D17").Value = "'0.005943"
$ws.Range("E17").Value = "'1.35%"
$ws.Range("D18").Value = "'44.09"
$ws.Range("E18").Value = "'-1.19%"
$ws.Range("D19").Value = "'3.359"
$ws.Range("E19").Value = "'-6.50%"
$ws.Range("D20").Value = "'0.3475"
$ws.Range("E20").Value = "'-0.27%"
$ws.Range("D21").Value = "'0.1407"
$ws.Range("E21").Value = "'3.74%"
$ws.Range("E22").Value = "'-4.29%"
$ws.Range("D23").Value = "'0.001258"
$ws.Range("E23").Value = "'1.59%"
$ws.Range("D24").Value = "'0.004335"
$ws.Range("E24").Value = "'-2.80%"
$ws.Range("D25").Value = "'0.0001199"
$ws.Range("E25").Value = "'-0.33%"
$ws.Range("D26").Value = "'0.0003743"
$ws.Range("E26").Value = "'-6.44%"
$ws.Range("D38").Value = "'0.02580"
$ws.Range("E38").Value = "'-2.31%"
$ws.Range("D39").Value = "'0.05628"
$ws.Range("E39").Value = "'2.59%"
$ws.Range("D40").Value = "'0.007566"
$ws.Range("E40").Value = "'-1.70%"
$ws.Range("D41").Value = "'0.1398"
$ws.Range("E41").Value = "'0.21%"
$ws.Range("D42").Value = "'0.007509"
$ws.Range("E42").Value = "'2.08%"
$ws.Range("D43").Value = "'0.002019"
$ws.Range("E43").Value = "'-1.62%"
$ws.Range("D44").Value = "'0.008876"
$ws.Range("E44").Value = "'1.89%"
$ws.Range("D45").Value = "'0.00007143"
$ws.Range("E45").Value = "'0.50%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.32%"
$ws.Range("E47").Value = "'54.99%"
$ws.Range("D48").Value = "'0.003600"
$ws.Range("E48").Value = "'1.78%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.32%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.32%"
